$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenes")

$ws.Range("C2").Value = "{name}, you live in {location:hometown:name}, a small village {location:hometown:feature:relativeposition}. You work {industry:hometown:workplace} every day, bringing back enough {industry:hometown:goods} {industry:hometown:purpose} your small village.`n`nOne day, a messenger rides into town. `"Citizens of {location:hometown:name}, the {character:antagonist:baron:baron} {character:baron:name} lays claim to your city. You are now under {character:baron:possPronoun} rule and will pay taxes and fealty to {character:baron:objPronoun}.`"`n`nWhat do you do?"

$ws.Range("G2").Value = "{industry:hometown:goodday} {|ROC|}"

$ws.Range("D25").Value = "Go back to your old job, {industry:hometown:workGer}"

$ws.Range("F25").Value = "{industry:hometown:gooddayfinal} {location:hometown:name}'s freedom allows things to go back to the way that they were. But you know that, if ever another {character:baron:baron} tried to take control of {location:hometown:name} again, you and your neighbors would be ready for them."

$ws.Rows.Item(2).RowHeight = 124.6
$ws.Rows.Item(13).RowHeight = 57.45
$ws.Rows.Item(14).RowHeight = 68.65
$ws.Rows.Item(25).RowHeight = 46.25

$ws.Range("D25").Select()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 2
